# Update column G ("K" - strikeouts) values per regenerated save_data.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- this script writes the recalculated K column values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 2
    6  = 2
    7  = 2
    8  = 2
    9  = 2
    10 = 1
    11 = 2
    12 = 2
    13 = 2
    14 = 1
    15 = 1
    17 = 1
    18 = 2
    19 = 2
    20 = 2
    21 = 2
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 0
    28 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
